$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 489.7143
$ws.Range("I33").Value = 159.16667
$ws.Range("K33").Value = 159.16667
$ws.Range("M33").Value = 69.83332999999999

$ws.Range("H41").Value = 1439.7693
$ws.Range("I41").Value = 450
$ws.Range("J41").Value = 1619.7273
$ws.Range("K41").Value = 450
$ws.Range("L41").Value = 1619.7273
$ws.Range("M41").Value = -10
$ws.Range("N41").Value = -2499.7273

$ws.Range("H64").Value = 46196.2
$ws.Range("I64").Value = 70660.336
$ws.Range("K64").Value = 70660.336
$ws.Range("M64").Value = -70412.336

$ws.Range("H67").Value = 46196.2
$ws.Range("I67").Value = 70660.336
$ws.Range("K67").Value = 70660.336
$ws.Range("M67").Value = -69802.336

$ws.Range("H107").Value = 1700.2632
$ws.Range("I107").Value = 1312.1177
$ws.Range("K107").Value = 1312.1177
$ws.Range("M107").Value = 607.8823

$ws.Range("H112").Value = 5077.25
$ws.Range("J112").Value = 5577.643
$ws.Range("L112").Value = 16732.929
$ws.Range("N112").Value = -18948.929

$ws.Range("H131").Value = 3983.5557
$ws.Range("I131").Value = 3113.6667
$ws.Range("J131").Value = 8333
$ws.Range("K131").Value = 9341.000100000001
$ws.Range("L131").Value = 24999
$ws.Range("M131").Value = -4301.000100000001
$ws.Range("N131").Value = -35079

$ws.Range("H135").Value = 982.04
$ws.Range("I135").Value = 473.85
$ws.Range("K135").Value = 4264.650000000001
$ws.Range("M135").Value = -1729.650000000001

$ws.Range("H138").Value = 33874.53
$ws.Range("I138").Value = 2208.682
$ws.Range("K138").Value = 6626.045999999999
$ws.Range("M138").Value = -1486.045999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14200.156
$ws.Range("I32").Value = 15037.208
$ws.Range("K32").Value = 15037.208
$ws.Range("M32").Value = -14750.208

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 4333.3335
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 4333.3335
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -5705.3335

$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 4333.3335
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 21666.6675
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -28530.6675

$ws.Range("H88").Value = 4311.25
$ws.Range("I88").Value = 1318.75
$ws.Range("J88").Value = 5308.75
$ws.Range("K88").Value = 1318.75
$ws.Range("L88").Value = 5308.75
$ws.Range("M88").Value = -912.75
$ws.Range("N88").Value = -6120.75

$ws.Range("H91").Value = 4311.25
$ws.Range("I91").Value = 1318.75
$ws.Range("J91").Value = 5308.75
$ws.Range("K91").Value = 1318.75
$ws.Range("L91").Value = 5308.75
$ws.Range("M91").Value = 85.25
$ws.Range("N91").Value = -8116.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1176.9778
$ws.Range("I86").Value = 1175.0625
$ws.Range("J86").Value = 1181.6923
$ws.Range("K86").Value = 1175.0625
$ws.Range("L86").Value = 1181.6923
$ws.Range("M86").Value = -52.0625
$ws.Range("N86").Value = -3427.6923

$ws.Range("H89").Value = 1176.9778
$ws.Range("I89").Value = 1175.0625
$ws.Range("J89").Value = 1181.6923
$ws.Range("K89").Value = 5875.3125
$ws.Range("L89").Value = 5908.461499999999
$ws.Range("M89").Value = -259.3125
$ws.Range("N89").Value = -17140.4615

$ws.Range("H132").Value = 95979
$ws.Range("J132").Value = 95979
$ws.Range("L132").Value = 95979
$ws.Range("N132").Value = -106099

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 23000
$ws.Range("J4").Value = 80250
$ws.Range("L4").Value = 80250
$ws.Range("N4").Value = -80474

$ws.Range("H6").Value = 165592.94
$ws.Range("I6").Value = 1714.4286
$ws.Range("K6").Value = 1714.4286
$ws.Range("M6").Value = -1601.4286

$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5340

$ws.Range("H17").Value = 16251.75
$ws.Range("I17").Value = 16251.75
$ws.Range("K17").Value = 16251.75
$ws.Range("M17").Value = -16077.75

$ws.Range("H31").Value = 3335362.5
$ws.Range("I31").Value = 4348901
$ws.Range("K31").Value = 4348901
$ws.Range("M31").Value = -4348606

$ws.Range("H34").Value = 3335362.5
$ws.Range("I34").Value = 4348901
$ws.Range("K34").Value = 4348901
$ws.Range("M34").Value = -4348699

$ws.Range("H58").Value = 13586.407
$ws.Range("I58").Value = 1632.5
$ws.Range("J58").Value = 66183.60000000001
$ws.Range("K58").Value = 1632.5
$ws.Range("L58").Value = 66183.60000000001
$ws.Range("M58").Value = -1429.5
$ws.Range("N58").Value = -66589.60000000001

$ws.Range("H99").Value = 9087.1
$ws.Range("I99").Value = 7401.3335
$ws.Range("J99").Value = 11615.75
$ws.Range("K99").Value = 7401.3335
$ws.Range("L99").Value = 11615.75
$ws.Range("M99").Value = -5903.3335
$ws.Range("N99").Value = -14611.75

$ws.Range("H107").Value = 953.3158
$ws.Range("I107").Value = 889.61536
$ws.Range("J107").Value = 1091.3334
$ws.Range("K107").Value = 889.61536
$ws.Range("L107").Value = 1091.3334
$ws.Range("M107").Value = 1030.38464
$ws.Range("N107").Value = -4931.3334

$ws.Range("H126").Value = 9087.1
$ws.Range("I126").Value = 7401.3335
$ws.Range("J126").Value = 11615.75
$ws.Range("K126").Value = 22204.0005
$ws.Range("L126").Value = 34847.25
$ws.Range("M126").Value = -19734.0005
$ws.Range("N126").Value = -39787.25

$ws.Range("H134").Value = 1428.8948
$ws.Range("I134").Value = 1274.6
$ws.Range("K134").Value = 3823.8
$ws.Range("M134").Value = -1288.8

$ws.Range("H136").Value = 13586.407
$ws.Range("I136").Value = 1632.5
$ws.Range("J136").Value = 66183.60000000001
$ws.Range("K136").Value = 4897.5
$ws.Range("L136").Value = 198550.8
$ws.Range("M136").Value = -2347.5
$ws.Range("N136").Value = -203650.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 38587.3
$ws.Range("I48").Value = 25990
$ws.Range("J48").Value = 39987
$ws.Range("K48").Value = 25990
$ws.Range("L48").Value = 39987
$ws.Range("M48").Value = -25505
$ws.Range("N48").Value = -40957

$ws.Range("H97").Value = 1452.9
$ws.Range("J97").Value = 2369.6365
$ws.Range("L97").Value = 2369.6365
$ws.Range("N97").Value = -3361.6365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1018.7368
$ws.Range("I16").Value = 1275.0714
$ws.Range("K16").Value = 1275.0714
$ws.Range("M16").Value = -1105.0714

$ws.Range("H40").Value = 1718.8611
$ws.Range("I40").Value = 1625.1143
$ws.Range("K40").Value = 1625.1143
$ws.Range("M40").Value = -1489.1143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3400
$ws.Range("I62").Value = 3400
$ws.Range("K62").Value = 3400
$ws.Range("M62").Value = -2776

$ws.Range("H65").Value = 3400
$ws.Range("I65").Value = 3400
$ws.Range("K65").Value = 17000
$ws.Range("M65").Value = -13880

$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180

$ws.Range("H116").Value = 106000
$ws.Range("J116").Value = 106000
$ws.Range("L116").Value = 106000
$ws.Range("N116").Value = -115178
